# Timesheet: finish off the still-open clock-in on row 13 (add its
# clock-out + duration) and append two more clock-in/out entries for
# 2026-02-03 (row 14 complete, row 15 clocked-in only, still open).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Write the new text values -------------------------------------
# A leading apostrophe forces these to be stored as literal text (not
# auto-parsed into a date/time serial number) while keeping today's
# General number format, matching how the rest of the sheet stores its
# dates/times/durations as plain strings.

# Close out row 13 (previously blank clock-out/duration).
$ws.Range("C13").Value = "'08:44:48"
$ws.Range("D13").Value = "'0.76 Hours"

# New row 14 - a complete clock-in/out entry.
$ws.Range("A14").Value = "'2026-02-03"
$ws.Range("B14").Value = "'12:52:45"
$ws.Range("C14").Value = "'13:01:55"
$ws.Range("D14").Value = "'0.15 Hours"

# New row 15 - only clocked in so far; clock-out/duration stay empty.
$ws.Range("A15").Value = "'2026-02-03"
$ws.Range("B15").Value = "'16:15:23"

# --- 2) Match the formatting of the rest of the table ------------------
# Re-apply the existing data-row style (copied from an already-styled
# cell) on top of the values so the new/edited cells line up with the
# rest of the sheet instead of picking up a "quoted text" style.

$ws.Range("B13").Copy()
$ws.Range("C13:D13").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A13:D13").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A15:D15").PasteSpecial(-4122)  # xlPasteFormats (covers blank C15/D15 too)

$excel.CutCopyMode = 0
